$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 64-66 were ROPME / ROPME.n, should become PERSGA / PERSGA.n
$ws.Range("A64").Value = "PERSGA"
$ws.Range("B64").Value = "PERSGA.1"

$ws.Range("A65").Value = "PERSGA"
$ws.Range("B65").Value = "PERSGA.2"

$ws.Range("A66").Value = "PERSGA"
$ws.Range("B66").Value = "PERSGA.3"

# Rows 67-69 were PERSGA / PERSGA.n, should become ROPME / ROPME.n
$ws.Range("A67").Value = "ROPME"
$ws.Range("B67").Value = "ROPME.1"

$ws.Range("A68").Value = "ROPME"
$ws.Range("B68").Value = "ROPME.2"

$ws.Range("A69").Value = "ROPME"
$ws.Range("B69").Value = "ROPME.3"

# Update view state: scroll the frozen pane down and change the selection
# to the newly-relevant rows (67-69), matching where the author was
# reviewing the ROPME / PERSGA mix-up.
$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("A67:XFD69").Select()
